# The commit adds one new weekly price record for "Poroto granado" at
# "Feria Lagunitas de Puerto Montt". In the source data the rows are kept in
# reverse-chronological order, so the new record is inserted as the sheet's
# 6th data row (worksheet row 7, just under the header row and the five most
# recent existing records), pushing every following row down by one.
#
# Net effect mirrors the diff: dimension grows from A1:R57 to A1:R58, rows
# 7-57 become rows 8-58 (content unchanged, just shifted down), and the new
# row 7 holds the values below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7:57 down to 8:58, leaving a blank row 7 (inherits formatting,
# including the date style already used by column D, from the row above).
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new market record.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C7").Value = 'Los Lagos'
$ws.Range("D7").Value = 44950
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112030
$ws.Range("G7").Value = 'Poroto granado'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 48000
$ws.Range("L7").Value = 48000
$ws.Range("M7").Value = 48000
$ws.Range("N7").Value = '$/saco 25 kilos'
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 1920
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = 'Hortaliza'
